# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet — the third sheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$legislatorName = "謝國樑"
$legislatorId = 1387
$reportDate = "2012-05-01"

# --- Header row (row 1): date / legislator_name / legislator_id ---
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Match the bold/border/centered header formatting used by the existing
# header cells (e.g. G1) by copying its format onto the new header cells.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-14) ---
for ($r = 2; $r -le 14; $r++) {
    # Write the date as a quoted-text formula first so Excel stores it as
    # literal text ("2012-05-01") instead of auto-converting it to a date
    # serial number, then collapse the formula down to its plain value.
    $ws.Cells.Item($r, 8).Formula = '="' + $reportDate + '"'
    $ws.Cells.Item($r, 8).Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4163)

    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
$excel.CutCopyMode = $false

# Match the plain data-row formatting used by the other data columns (e.g.
# column C) by copying its format onto the new data cells, row by row.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Copy()
    $ws.Range("H" + $r + ":J" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
